$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get a plain-decimal string (e.g. "306.07").
# Temporarily format them as Text so the literal is stored verbatim
# (matching the source, which keeps the whole Price column as text)
# instead of Excel silently converting the input to a Number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.461.69"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.237.05"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +1.49%  "
$ws.Range("D5").Value = "306.07"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "94.09"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "34.65"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "7.19"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D14").Value = "0.833"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "2.209.65"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "13.54"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "44.119.51"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "0.0₃0952"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "6.31"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").Value = "11.91"
$ws.Range("D21").Value = "65.44"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").Value = "236.88"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  +3.19%  "
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").Value = "38.40"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").Value = "5.93"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "19.92"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "153.09"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "0.0792"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("D34").Value = "3.05"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").Value = "14.90"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "3.75"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "1.787.43"
$ws.Range("E43").Value = "  +3.49%  "
$ws.Range("D44").Value = "0.191"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").Value = "78.45"
$ws.Range("E45").Value = "  -8.22%  "
$ws.Range("D46").Value = "1.64"
$ws.Range("E46").Value = "  +9.31%  "
$ws.Range("D47").Value = "70.31"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").Value = "98.61"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").Value = "4.89"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "54.43"
$ws.Range("E51").Value = "  +0.58%  "

# Restore each cell's original (default) style now that the text
# value is safely stored, so no visible formatting changes remain.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
